# "Update for round 99"
#
# Adds survey rounds 97-99 (rows 110-112) to the UK sheet, and fixes up
# row 109 (round 96) whose received date / spss_name were corrected at
# the same time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Fix existing row 109 (survey_round 96): date_recieved moved from
#    2022-02-04 to 2022-01-31, and the spss_name corrected.
# ---------------------------------------------------------------------
$ws.Range("G109").Value = 44592
$ws.Range("H109").Value = "21-088043_PEW39_Final_ICUO"

# ---------------------------------------------------------------------
# 2. Append new rows 110-112 for survey rounds 97, 98 and 99.
#    Copy the formatting (date number format) of row 109 down first so
#    the new G cells pick up the same style as the rest of the column.
# ---------------------------------------------------------------------
$ws.Range("G109").Copy() | Out-Null
$ws.Range("G110:G112").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 110 - survey_round 97, panel F, wave 39
$ws.Range("A110").Value = 8
$ws.Range("B110").Value = 0
$ws.Range("C110").Value = "uk"
$ws.Range("D110").Value = 97
$ws.Range("E110").Value = "F"
$ws.Range("F110").Value = 39
$ws.Range("G110").Value = 44596
$ws.Range("H110").Value = "21-088071_PFW39_Final_ICUO"
$ws.Range("J110").Value = 1
$ws.Range("K110").Value = 1

# Row 111 - survey_round 98, panel E, wave 40
$ws.Range("A111").Value = 8
$ws.Range("B111").Value = 0
$ws.Range("C111").Value = "uk"
$ws.Range("D111").Value = 98
$ws.Range("E111").Value = "E"
$ws.Range("F111").Value = 40
$ws.Range("G111").Value = 44602
$ws.Range("H111").Value = "21-088043_PEW40_Final_ICUO"
$ws.Range("J111").Value = 1
$ws.Range("K111").Value = 1

# Row 112 - survey_round 99, panel F, wave 40
$ws.Range("A112").Value = 8
$ws.Range("B112").Value = 0
$ws.Range("C112").Value = "uk"
$ws.Range("D112").Value = 99
$ws.Range("E112").Value = "F"
$ws.Range("F112").Value = 40
$ws.Range("G112").Value = 44610
$ws.Range("H112").Value = "21-088071_PFW40_Final_ICUO"
$ws.Range("J112").Value = 1
$ws.Range("K112").Value = 1

# ---------------------------------------------------------------------
# 3. Extend the "spss_name" helper formula (column I) from I83:I109 down
#    to I83:I112, matching how the existing rows build the file name.
# ---------------------------------------------------------------------
for ($r = 110; $r -le 112; $r++) {
    $ws.Range("I$r").Formula = "=C$r&""_""&""sr""&TEXT(D$r,""00"")&""_""&YEAR(G$r)&TEXT(G$r,""MM"")&TEXT(G$r,""DD"")&""_p""&E$r&""_wv""&TEXT(F$r,""00"")&"""""
}

# ---------------------------------------------------------------------
# 4. Update the view so the newly added rows are visible / selected,
#    matching the author's final cursor position.
# ---------------------------------------------------------------------
$ws.Range("A88").Select() | Out-Null
$ws.Range("I112").Select() | Out-Null
